$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.022018107476085
$ws.Range("D2").Value = 1.032962751587116
$ws.Range("E2").Value = 1.022823282317981
$ws.Range("F2").Value = 1.044095121657267
$ws.Range("I2").Value = 1.033642522701657
$ws.Range("J2").Value = 1.027206307857129
$ws.Range("K2").Value = 1.03576646986476
$ws.Range("L2").Value = 1.025656529141747
$ws.Range("M2").Value = 1.046867136853178
$ws.Range("N2").Value = 1.013139812424245
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.022853609180232
$ws.Range("D3").Value = 1.033600382175798
$ws.Range("E3").Value = 1.023528897289025
$ws.Range("F3").Value = 1.044897680723172
$ws.Range("I3").Value = 1.033791518447116
$ws.Range("J3").Value = 1.027680270806949
$ws.Range("K3").Value = 1.036213439642323
$ws.Range("L3").Value = 1.026169168396784
$ws.Range("M3").Value = 1.047480886780453
$ws.Range("N3").Value = 1.013296993306183
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.02339474965871
$ws.Range("D4").Value = 1.034013043331188
$ws.Range("E4").Value = 1.023986325646339
$ws.Range("F4").Value = 1.045417349796081
$ws.Range("I4").Value = 1.033886289418028
$ws.Range("J4").Value = 1.02798683194037
$ws.Range("K4").Value = 1.03650204324569
$ws.Range("L4").Value = 1.026501056228335
$ws.Range("M4").Value = 1.047877729965687
$ws.Range("N4").Value = 1.013398634522934
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.023622366801488
$ws.Range("D5").Value = 1.034186541319721
$ws.Range("E5").Value = 1.02417882980115
$ws.Range("F5").Value = 1.04563590272557
$ws.Range("I5").Value = 1.03392573777882
$ws.Range("J5").Value = 1.02811567912987
$ws.Range("K5").Value = 1.036623223349677
$ws.Range("L5").Value = 1.026640622621399
$ws.Range("M5").Value = 1.048044490670786
$ws.Range("N5").Value = 1.013441348396689
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.023660591810276
$ws.Range("D6").Value = 1.034215673234101
$ws.Range("E6").Value = 1.024211163845906
$ws.Range("F6").Value = 1.045672603592923
$ws.Range("I6").Value = 1.033932338238871
$ws.Range("J6").Value = 1.02813731127818
$ws.Range("K6").Value = 1.036643561251601
$ws.Range("L6").Value = 1.026664058782392
$ws.Range("M6").Value = 1.048072486243252
$ws.Range("N6").Value = 1.013448519287045
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.023397790613988
$ws.Range("D7").Value = 1.034015361562306
$ws.Range("E7").Value = 1.023988897107691
$ws.Range("F7").Value = 1.045420269783355
$ws.Range("I7").Value = 1.03388681807573
$ws.Range("J7").Value = 1.027988553726685
$ws.Range("K7").Value = 1.036503663046874
$ws.Range("L7").Value = 1.026502920962928
$ws.Range("M7").Value = 1.047879958516204
$ws.Range("N7").Value = 1.013399205331617
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.022300361586107
$ws.Range("D8").Value = 1.033178225888122
$ws.Range("E8").Value = 1.023061571471347
$ws.Range("F8").Value = 1.044366274537532
$ws.Range("I8").Value = 1.033693215375671
$ws.Range("J8").Value = 1.027366510880474
$ws.Range("K8").Value = 1.035917652018142
$ws.Range("L8").Value = 1.025829740719718
$ws.Range("M8").Value = 1.047074616337178
$ws.Range("N8").Value = 1.013192945691773
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.020370566785651
$ws.Range("D9").Value = 1.031703716737339
$ws.Range("E9").Value = 1.021434081883412
$ws.Range("F9").Value = 1.042511845838738
$ws.Range("I9").Value = 1.033339548458048
$ws.Range("J9").Value = 1.026269494337774
$ws.Range("K9").Value = 1.034880373534089
$ws.Range("L9").Value = 1.024644919369313
$ws.Range("M9").Value = 1.04565332644607
$ws.Range("N9").Value = 1.012829009606376
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.019086831862168
$ws.Range("D10").Value = 1.030721246521219
$ws.Range("E10").Value = 1.020353618115283
$ws.Range("F10").Value = 1.041277595279082
$ws.Range("I10").Value = 1.033095410801104
$ws.Range("J10").Value = 1.025537622591411
$ws.Range("K10").Value = 1.034185812090936
$ws.Range("L10").Value = 1.023856066897933
$ws.Range("M10").Value = 1.044704442998522
$ws.Range("N10").Value = 1.012586089628553
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.018531643317306
$ws.Range("D11").Value = 1.030295976675152
$ws.Range("E11").Value = 1.019886862334504
$ws.Range("F11").Value = 1.040743659524208
$ws.Range("I11").Value = 1.032987724995809
$ws.Range("J11").Value = 1.025220604928175
$ws.Range("K11").Value = 1.033884355602827
$ws.Range("L11").Value = 1.023514746163959
$ws.Range("M11").Value = 1.044293266169403
$ws.Range("N11").Value = 1.012480838073653
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.018325524894093
$ws.Range("D12").Value = 1.030138036380223
$ws.Range("E12").Value = 1.019713654376653
$ws.Range("F12").Value = 1.040545409622482
$ws.Range("I12").Value = 1.032947430216045
$ws.Range("J12").Value = 1.025102834847986
$ws.Range("K12").Value = 1.033772276479584
$ws.Range("L12").Value = 1.023388004526771
$ws.Range("M12").Value = 1.044140492868874
$ws.Range("N12").Value = 1.012441733557516
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.018369733294511
$ws.Range("D13").Value = 1.030171913999587
$ws.Range("E13").Value = 1.019750800533847
$ws.Range("F13").Value = 1.040587931354314
$ws.Range("I13").Value = 1.03295608693991
$ws.Range("J13").Value = 1.025128097615644
$ws.Range("K13").Value = 1.033796322549775
$ws.Range("L13").Value = 1.023415189202316
$ws.Range("M13").Value = 1.044173265219654
$ws.Range("N13").Value = 1.012450122027176
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.018514603387217
$ws.Range("D14").Value = 1.030282920783439
$ws.Range("E14").Value = 1.019872541515324
$ws.Range("F14").Value = 1.040727270530357
$ws.Range("I14").Value = 1.032984400239292
$ws.Range("J14").Value = 1.025210870327566
$ws.Range("K14").Value = 1.033875093236394
$ws.Range("L14").Value = 1.023504268835086
$ws.Range("M14").Value = 1.044280638766493
$ws.Range("N14").Value = 1.012477605871336
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.018603876371349
$ws.Range("D15").Value = 1.030351318944895
$ws.Range("E15").Value = 1.019947572193027
$ws.Range("F15").Value = 1.040813132341204
$ws.Range("I15").Value = 1.033001805870272
$ws.Range("J15").Value = 1.025261867297582
$ws.Range("K15").Value = 1.033923612616923
$ws.Range("L15").Value = 1.023559159089153
$ws.Range("M15").Value = 1.044346789385373
$ws.Range("N15").Value = 1.012494538342453
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.019123692299719
$ws.Range("D16").Value = 1.030749473522916
$ws.Range("E16").Value = 1.020384618326683
$ws.Range("F16").Value = 1.041313041607671
$ws.Range("I16").Value = 1.033102516068538
$ws.Range("J16").Value = 1.025558659751271
$ws.Range("K16").Value = 1.034205803995721
$ws.Range("L16").Value = 1.023878724776008
$ws.Range("M16").Value = 1.044731725195118
$ws.Range("N16").Value = 1.012593073484209
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.019449941511175
$ws.Range("D17").Value = 1.030999265859751
$ws.Range("E17").Value = 1.020659059520193
$ws.Range("F17").Value = 1.041626757574093
$ws.Range("I17").Value = 1.033165161420873
$ws.Range("J17").Value = 1.025744800634832
$ws.Range("K17").Value = 1.034382626971039
$ws.Range("L17").Value = 1.024079249828774
$ws.Range("M17").Value = 1.044973105103574
$ws.Range("N17").Value = 1.012654864732553
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.019640302384997
$ws.Range("D18").Value = 1.031144979487415
$ws.Range("E18").Value = 1.020819241679963
$ws.Range("F18").Value = 1.041809791172751
$ws.Range("I18").Value = 1.033201511009807
$ws.Range("J18").Value = 1.025853362535651
$ws.Range("K18").Value = 1.034485696529529
$ws.Range("L18").Value = 1.024196237460918
$ws.Range("M18").Value = 1.045113868467735
$ws.Range("N18").Value = 1.012690900162345
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.01970522153917
$ws.Range("D19").Value = 1.031194666383157
$ws.Range("E19").Value = 1.020873877454056
$ws.Range("F19").Value = 1.041872209043777
$ws.Range("I19").Value = 1.033213872960758
$ws.Range("J19").Value = 1.02589037744598
$ws.Range("K19").Value = 1.034520828958514
$ws.Range("L19").Value = 1.024236131378176
$ws.Range("M19").Value = 1.045161860087988
$ws.Range("N19").Value = 1.012703186220092
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.019414931300011
$ws.Range("D20").Value = 1.030972464045957
$ws.Range("E20").Value = 1.020629603677157
$ws.Range("F20").Value = 1.041593093800504
$ws.Range("I20").Value = 1.03315845985349
$ws.Range("J20").Value = 1.025724830600679
$ws.Range("K20").Value = 1.034363662592603
$ws.Range("L20").Value = 1.024057732817193
$ws.Range("M20").Value = 1.044947210357499
$ws.Range("N20").Value = 1.012648235770953
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.018471939890754
$ws.Range("D21").Value = 1.030250231378539
$ws.Range("E21").Value = 1.019836687235623
$ws.Range("F21").Value = 1.040686236456998
$ws.Range("I21").Value = 1.032976070827855
$ws.Range("J21").Value = 1.025186496247489
$ws.Range("K21").Value = 1.033851900112973
$ws.Range("L21").Value = 1.023478035985993
$ws.Range("M21").Value = 1.044249021131429
$ws.Range("N21").Value = 1.012469512820162
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.017879642341013
$ws.Range("D22").Value = 1.029796274018638
$ws.Range("E22").Value = 1.019339110074051
$ws.Range("F22").Value = 1.040116509301198
$ws.Range("I22").Value = 1.032859686504451
$ws.Range("J22").Value = 1.024847934418463
$ws.Range("K22").Value = 1.033529529317526
$ws.Range("L22").Value = 1.023113790155714
$ws.Range("M22").Value = 1.043809788621867
$ws.Range("N22").Value = 1.012357088432243
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.018193573210184
$ws.Range("D23").Value = 1.030036911600884
$ws.Range("E23").Value = 1.019602793449887
$ws.Range("F23").Value = 1.04041848910064
$ws.Range("I23").Value = 1.032921545745847
$ws.Range("J23").Value = 1.025027420522383
$ws.Range("K23").Value = 1.033700481086326
$ws.Range("L23").Value = 1.023306861256543
$ws.Range("M23").Value = 1.04404265750316
$ws.Range("N23").Value = 1.012416691711467
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.019430750706071
$ws.Range("D24").Value = 1.030984574588962
$ws.Range("E24").Value = 1.020642913179725
$ws.Range("F24").Value = 1.041608304861497
$ws.Range("I24").Value = 1.033161488591677
$ws.Range("J24").Value = 1.025733854233426
$ws.Range("K24").Value = 1.03437223199008
$ws.Range("L24").Value = 1.024067455352196
$ws.Range("M24").Value = 1.044958911169832
$ws.Range("N24").Value = 1.012651231133064
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.020868979343792
$ws.Range("D25").Value = 1.032084826398545
$ws.Range("E25").Value = 1.021854036488782
$ws.Range("F25").Value = 1.042990910295483
$ws.Range("I25").Value = 1.033432457242266
$ws.Range("J25").Value = 1.026553197924306
$ws.Range("K25").Value = 1.035149077374232
$ws.Range("L25").Value = 1.024951048782739
$ws.Range("M25").Value = 1.046021010763525
$ws.Range("N25").Value = 1.012923149936152
